# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on the per-job Leve sheets with refreshed Universalis
# pull data. Generated from the upstream diff; applied cell-by-cell so any
# cell that should no longer carry a profit value (K/L price == 0) is cleared
# rather than zeroed, matching how this sheet represents "no data".

$wb = $excel.ActiveWorkbook

# Map of SheetName -> RowNumber -> @{ Column = NewValue (or $null to clear) }
$updates = @{
    "ALC" = @{
        52 = @{ "H" = 0; "I" = 0; "K" = 0; "M" = $null }
        64 = @{ "H" = 6498.1665; "I" = 4999.5; "J" = 7247.5; "K" = 4999.5; "L" = 7247.5; "M" = -4751.5; "N" = -7743.5 }
        67 = @{ "H" = 6498.1665; "I" = 4999.5; "J" = 7247.5; "K" = 4999.5; "L" = 7247.5; "M" = -4141.5; "N" = -8963.5 }
        70 = @{ "H" = 2166.6667; "J" = 2333.3333; "L" = 6999.999899999999; "N" = -7539.999899999999 }
        73 = @{ "H" = 2166.6667; "J" = 2333.3333; "L" = 6999.999899999999; "N" = -8871.999899999999 }
        74 = @{ "H" = 7835.8184; "I" = 5998.5; "J" = 8885.714; "K" = 5998.5; "L" = 8885.714; "M" = -5062.5; "N" = -10757.714 }
        77 = @{ "H" = 7835.8184; "I" = 5998.5; "J" = 8885.714; "K" = 29992.5; "L" = 44428.57; "M" = -25312.5; "N" = -53788.57 }
        92 = @{ "H" = 1390.625; "I" = 1517.8572; "K" = 1517.8572; "M" = -269.8571999999999 }
        99 = @{ "H" = 276.5; "I" = 276.5; "J" = 0; "K" = 829.5; "L" = 0; "M" = 668.5; "N" = $null }
    }
    "ARM" = @{
        32 = @{ "H" = 3831.48; "I" = 3599.4783; "K" = 3599.4783; "M" = -3312.4783 }
        44 = @{ "H" = 63049; "J" = 63049; "L" = 63049; "N" = -64025 }
        63 = @{ "H" = 11688.125; "I" = 9083.166999999999; "J" = 19503; "K" = 9083.166999999999; "L" = 19503; "M" = -8397.166999999999; "N" = -20875 }
        66 = @{ "H" = 11688.125; "I" = 9083.166999999999; "J" = 19503; "K" = 45415.835; "L" = 97515; "M" = -41983.835; "N" = -104379 }
    }
    "BSM" = @{
        22 = @{ "H" = 229.6; "I" = 229.6; "K" = 229.6; "M" = -56.59999999999999 }
        80 = @{ "H" = 1343.3334; "I" = 85.25; "K" = 85.25; "M" = 912.75 }
        83 = @{ "H" = 1343.3334; "I" = 85.25; "K" = 426.25; "M" = 4565.75 }
        86 = @{ "H" = 1200; "I" = 0; "K" = 0; "M" = $null }
        89 = @{ "H" = 1200; "I" = 0; "K" = 0; "M" = $null }
        94 = @{ "H" = 729.375; "I" = 672.6667; "K" = 672.6667; "M" = -221.6667 }
        107 = @{ "H" = 969; "I" = 994; "J" = 956.5; "K" = 994; "L" = 956.5; "M" = 926; "N" = -4796.5 }
    }
    "CRP" = @{
        22 = @{ "H" = 50; "I" = 50; "K" = 50; "M" = 300 }
        62 = @{ "H" = 3666.6667; "J" = 3750; "L" = 3750; "N" = -4998 }
        65 = @{ "H" = 3666.6667; "J" = 3750; "L" = 18750; "N" = -24990 }
    }
    "CUL" = @{
        36 = @{ "H" = 550; "I" = 550; "J" = 0; "K" = 1650; "L" = 0; "M" = -1481; "N" = $null }
        46 = @{ "H" = 4056; "J" = 4056; "L" = 12168; "N" = -12350 }
        60 = @{ "H" = 0; "I" = 0; "K" = 0; "M" = $null }
        131 = @{ "H" = 2030.625; "J" = 2721.6667; "L" = 8165.000100000001; "N" = -18245.0001 }
    }
    "GSM" = @{
        122 = @{ "H" = 2720.4; "I" = 2755.2222; "K" = 8265.6666; "M" = -5815.6666 }
    }
    "LTW" = @{
        22 = @{ "H" = 1424.1666; "I" = 849.1667; "J" = 1999.1666; "K" = 849.1667; "L" = 1999.1666; "M" = -554.1667; "N" = -2589.1666 }
        27 = @{ "H" = 1424.1666; "I" = 849.1667; "J" = 1999.1666; "K" = 849.1667; "L" = 1999.1666; "M" = -742.1667; "N" = -2213.1666 }
        46 = @{ "H" = 2100; "I" = 2000; "J" = 2600; "K" = 2000; "L" = 2600; "M" = -1812; "N" = -2976 }
        55 = @{ "H" = 2113.5557; "I" = 2432.7144; "K" = 2432.7144; "M" = -2259.7144 }
        68 = @{ "H" = 5000; "I" = 5000; "K" = 5000; "M" = -4251 }
        71 = @{ "H" = 5000; "I" = 5000; "K" = 25000; "M" = -21256 }
        93 = @{ "H" = 1745.5; "I" = 1843; "K" = 1843; "M" = -595 }
    }
    "WVR" = @{
        100 = @{ "H" = 1262.5; "J" = 1624; "L" = 3248; "N" = -4330 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Sheets($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $cols = $rows[$rowNum]
        foreach ($colLetter in $cols.Keys) {
            $cellRef = "$colLetter$rowNum"
            $newValue = $cols[$colLetter]
            if ($null -eq $newValue) {
                $ws.Range($cellRef).ClearContents()
            } else {
                $ws.Range($cellRef).Value = $newValue
            }
        }
    }
}